# Add new power plant types to the Electricity Source subscript
# (issues #280 and #99)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoFCRfP")

$newPlants = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
for ($i = 0; $i -lt $newPlants.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newPlants[$i]
    $ws.Cells.Item($row, 2).Value = 1
}

# Update the selection on the SoFCRfP sheet to the cell below the new data
$ws.Range("A25").Select()

# Give the sheet a proper page setup (portrait orientation)
$ws.PageSetup.Orientation = 1

# Make "About" the active/selected sheet, as it is in the final workbook
$about = $wb.Worksheets.Item("About")
$about.Activate()
